$d = $word.ActiveDocument
$d.Content.Find.Execute("Ich brauche ein Dokument ich verändere", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
$d.Content.Find.Execute(" was", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
